$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that Excel would otherwise auto-coerce to a number
# (e.g. "247.96") while preserving it as plain text and keeping the cell's
# original (default/General) style -- matches how these price cells were
# already stored as text in the source workbook.
function Set-TextValue {
    param($cellRef, $val)
    $r = $ws.Range($cellRef)
    $r.Value = "'" + $val
    $r.Style = "Normal"
}

# Helper for plain non-numeric-looking text - no special handling needed.
function Set-PlainValue {
    param($cellRef, $val)
    $ws.Range($cellRef).Value = $val
}

# --- Price (column D) updates ---
Set-TextValue "D2"  "247.96"
Set-TextValue "D3"  "22.46"
Set-TextValue "D4"  "5.243"
Set-TextValue "D5"  "0.05688"
Set-TextValue "D7"  "6.314"
Set-TextValue "D8"  "0.8067"
Set-TextValue "D9"  "0.8778"
Set-TextValue "D11" "0.07396"
Set-TextValue "D12" "0.03051"
Set-TextValue "D14" "0.09395"
Set-TextValue "D15" "3.871"
Set-TextValue "D16" "0.001576"
Set-TextValue "D17" "0.04782"
Set-TextValue "D18" "0.0005810"
Set-TextValue "D19" "0.006400"
Set-TextValue "D20" "0.005040"
Set-TextValue "D21" "0.0009972"
Set-TextValue "D23" "3.691"
Set-TextValue "D24" "2.199"
Set-TextValue "D25" "0.3281"
Set-TextValue "D26" "0.1357"
Set-TextValue "D27" "0.0004750"
Set-TextValue "D45" "0.00005590"
Set-TextValue "D47" "0.4500"
Set-TextValue "D48" "0.1962"
Set-TextValue "D49" "0.00002100"

# --- Rows 41-43: coin ranking rotated by one position ---
# Row 41: BKEXToken -> KickToken
Set-PlainValue "B41" "KickToken"
Set-PlainValue "C41" "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue  "D41" "0.006831"
Set-PlainValue "E41" "40KickTokenKICK"

# Row 42: CEJI -> BKEXToken
Set-PlainValue "B42" "BKEXToken"
Set-PlainValue "C42" "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue  "D42" "0.1069"
Set-PlainValue "E42" "41BKEXTokenBKK"

# Row 43: KickToken -> CEJI
Set-PlainValue "B43" "CEJI"
Set-PlainValue "C43" "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue  "D43" "0.002730"
Set-PlainValue "E43" "42CEJICEJI"

# --- Row 48: "Worst in 24h" marker moved onto BOLO ---
Set-PlainValue "E48" "47BOLOBOLOWorstin24h"
